# Add a "Stage" column (AP) to the Observed sheet and stamp every
# observation row with the HarvestRipe stage name.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Observed")

# New header in AP1
$ws.Range("AP1").Value = "Stage"

# Data rows (2 through 83) all record the HarvestRipe stage
$ws.Range("AP2:AP83").Value = "HarvestRipe"

# Leave the sheet scrolled/selected near the newly added column, as in the
# authored workbook (cell AJ9 selected).
$ws.Range("AJ9").Select() | Out-Null
